$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.520102666666667
$ws.Range("H2").Value = 13.560308
$ws.Range("I2").Value = 0.9927775608668273
$ws.Range("J2").Value = 0.9927775608668273
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.166981632712
$ws.Range("R2").Value = 1.502834694408
$ws.Range("S2").Value = 0.02083872791335888
$ws.Range("T2").Value = 0.02083872791335888

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.520102666666667
$ws.Range("H3").Value = 13.560308
$ws.Range("I3").Value = 0.9927775608668273
$ws.Range("J3").Value = 0.9927775608668273
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("Q3").Value = 4.262961559464445
$ws.Range("R3").Value = 38.36665403518001
$ws.Range("S3").Value = 0.5320027993498209
$ws.Range("T3").Value = 0.5320027993498209

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.520102666666667
$ws.Range("H4").Value = 13.560308
$ws.Range("I4").Value = 0.9927775608668273
$ws.Range("J4").Value = 0.9927775608668273
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("Q4").Value = 3.525226563032445
$ws.Range("R4").Value = 31.727039067292
$ws.Range("S4").Value = 0.4399360336036476
$ws.Range("T4").Value = 0.4399360336036477

$ws.Range("I5").Value = 0.007222439133172593
$ws.Range("J5").Value = 0.007222439133172593
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.036942
$ws.Range("N5").Value = 0.110826
$ws.Range("O5").Value = 0.02099032928903418
$ws.Range("P5").Value = 0.02099032928903418
$ws.Range("Q5").Value = 0.001214788414
$ws.Range("R5").Value = 0.010933095726
$ws.Range("S5").Value = 0.0001516013756752993
$ws.Range("T5").Value = 0.0001516013756752993

$ws.Range("I6").Value = 0.007222439133172593
$ws.Range("J6").Value = 0.007222439133172593
$ws.Range("O6").Value = 0.5358731102718634
$ws.Range("P6").Value = 0.5358731102718634
$ws.Range("S6").Value = 0.003870310922042418
$ws.Range("T6").Value = 0.003870310922042418

$ws.Range("I7").Value = 0.007222439133172593
$ws.Range("J7").Value = 0.007222439133172593
$ws.Range("O7").Value = 0.4431365604391025
$ws.Range("P7").Value = 0.4431365604391026
$ws.Range("S7").Value = 0.003200526835454876
$ws.Range("T7").Value = 0.003200526835454876

